# Update the Czech Republic M2 data sheet:
#  - Row 259 (last existing data row): C:F values corrected from 6243183470000 -> 6243090940000
#  - Append 3 new monthly rows (260-262) with the same layout/style as existing rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the existing last row (259) ---
$ws.Range("C259").Value = 6243090940000
$ws.Range("D259").Value = 6243090940000
$ws.Range("E259").Value = 6243090940000
$ws.Range("F259").Value = 6243090940000
# G259 stays 0 (unchanged)

# --- Append new rows 260-262 ---
# Copy the formatting of the last existing row down into the new rows first
$ws.Range("A259:G259").Copy()
$ws.Range("A260:G262").PasteSpecial(-4122)

$newRows = @(
    @{ Row = 260; Date = 45108.41666666666; O = 6355692770000; H = 6355692770000; L = 6355692770000; C = 6355692770000 },
    @{ Row = 261; Date = 45139.41666666666; O = 6337051350000; H = 6337051350000; L = 6337051350000; C = 6337051350000 },
    @{ Row = 262; Date = 45170.41666666666; O = 6359425540000; H = 6359425540000; L = 6359425540000; C = 6359425540000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "ECONOMICS:CZM2"
    $ws.Cells.Item($row, 3).Value = $r.O
    $ws.Cells.Item($row, 4).Value = $r.H
    $ws.Cells.Item($row, 5).Value = $r.L
    $ws.Cells.Item($row, 6).Value = $r.C
    $ws.Cells.Item($row, 7).Value = 0
}
